# Refresh the cryptocurrency price/volume figures (columns D and E) to match the
# latest snapshot from the scraping run. Column A (rank), B (coin name) and C (link)
# are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.409.87'
$ws.Range("D3").Value = '2.246.93'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '''308.03'
$ws.Range("E5").Value = '  +1.14%  '
$ws.Range("D6").Value = '''94.86'
$ws.Range("E6").Value = '  -1.48%  '
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("E9").Value = '  +0.65%  '
$ws.Range("D10").Value = '''35.02'
$ws.Range("E10").Value = '  +1.40%  '
$ws.Range("D11").Value = '''0.0811'
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("D12").Value = '''7.21'
$ws.Range("E12").Value = '  +0.74%  '
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").Value = '2.356.74'
$ws.Range("E14").Value = '  +5.03%  '
$ws.Range("D15").Value = '''0.841'
$ws.Range("E15").Value = '  +2.79%  '
$ws.Range("D16").Value = '''13.71'
$ws.Range("E16").Value = '  +0.93%  '
$ws.Range("D17").Value = '44.119.95'
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").Value = '0.0₃0967'
$ws.Range("E18").Value = '  +0.87%  '
$ws.Range("D19").Value = '''12.30'
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = '''6.40'
$ws.Range("E20").Value = '  +3.13%  '
$ws.Range("D21").Value = '''65.88'
$ws.Range("E21").Value = '  +1.92%  '
$ws.Range("E22").Value = '  +3.95%  '
$ws.Range("D23").Value = '''237.49'
$ws.Range("E23").Value = '  -0.34%  '
$ws.Range("E24").Value = '  +4.25%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '''38.49'
$ws.Range("E26").Value = '  +6.82%  '
$ws.Range("E27").Value = '  +4.82%  '
$ws.Range("D28").Value = '''9.88'
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").Value = '''5.97'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = '''20.10'
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("D31").Value = '''154.17'
$ws.Range("E31").Value = '  +0.53%  '
$ws.Range("D32").Value = '''0.0801'
$ws.Range("E32").Value = '  -0.50%  '
$ws.Range("E33").Value = '  -0.13%  '
$ws.Range("D34").Value = '''3.11'
$ws.Range("E34").Value = '  -7.77%  '
$ws.Range("E35").Value = '  +1.49%  '
$ws.Range("E36").Value = '  +3.17%  '
$ws.Range("D37").Value = '''1.81'
$ws.Range("E37").Value = '  +2.70%  '
$ws.Range("E38").Value = '  +6.21%  '
$ws.Range("D39").Value = '''14.80'
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("E41").Value = '  +1.53%  '
$ws.Range("E42").Value = '  +0.25%  '
$ws.Range("D43").Value = '1.741.24'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("E44").Value = '  +3.44%  '
$ws.Range("E45").Value = '  -5.98%  '
$ws.Range("D46").Value = '''100.09'
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").Value = '''4.96'
$ws.Range("E47").Value = '  -2.29%  '
$ws.Range("D48").Value = '''70.96'
$ws.Range("E48").Value = '  +4.59%  '
$ws.Range("D49").Value = '''56.17'
$ws.Range("E49").Value = '  +3.62%  '
$ws.Range("E50").Value = '  +6.19%  '
$ws.Range("E51").Value = '  -0.02%  '
